$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Taxonsorteringsordning (column B) values: 79243 -> 79244 and 91828 -> 91829
$rowsB79243 = 3,4,6,7,8,10,11,12,13,14,15,16
foreach ($r in $rowsB79243) {
    $ws.Range("B$r").Value = 79244
}

$rowsB91828 = 17,22,23
foreach ($r in $rowsB91828) {
    $ws.Range("B$r").Value = 91829
}

# Rows 6 and 7 swap their identity-specific data (A, Q, R, AC, AM, AO)
# New row 6 values (previously held by row 7)
$ws.Range("A6").Value = 131154323
$ws.Range("Q6").Value = 445033
$ws.Range("R6").Value = 7053229
$ws.Range("AC6").Value = ""
$ws.Range("AM6").Value = ""
$ws.Range("AO6").Value = "Picea abies"

# New row 7 values (previously held by row 6)
$ws.Range("A7").Value = 131154313
$ws.Range("Q7").Value = 445053
$ws.Range("R7").Value = 7053130
$ws.Range("AC7").Value = "På flera granar."
$ws.Range("AM7").Value = "Gren på levande träd"
$ws.Range("AO7").Value = "Branch on living tree # Picea abies"
